# "isActive" column needs to be in lowercase
#
# 1. The "Export Worksheet" sheet has an ISACTIVE column (column C) whose
#    141 data rows (C2:C142) all currently hold the text "TRUE". They need
#    to read "true" (lowercase) instead.
# 2. The "SQL" sheet contains the export query as a single text cell; the
#    `'TRUE' AS isActive` literal inside that SQL needs to become
#    `'true' AS isActive`.

$wb = $excel.ActiveWorkbook

# --- Export Worksheet: lowercase the ISACTIVE column values -----------------
$ws1 = $wb.Worksheets.Item("Export Worksheet")
$ws1.Activate()

$isActiveRange = $ws1.Range("C2:C142")
$rowCount = $isActiveRange.Rows.Count
$values = New-Object 'object[,]' $rowCount, 1
for ($i = 0; $i -lt $rowCount; $i++) {
    # A leading apostrophe forces Excel to store this as literal text
    # ("true") instead of auto-converting the boolean-looking word into a
    # real Boolean value.
    $values[$i, 0] = "'true"
}
$isActiveRange.Value = $values

# Drop the "quote prefix" cell style picked up above so the cells keep their
# original (default) formatting.
$isActiveRange.Style = "Normal"

# Restore the selection to match the edited workbook (active cell E14).
[void]$ws1.Range("E14").Select()

# --- SQL sheet: lowercase the isActive literal in the query text -----------
$ws2 = $wb.Worksheets.Item("SQL")
$sqlCell = $ws2.Range("A2")
$sqlText = $sqlCell.Text
$sqlText = $sqlText.Replace("'TRUE' AS isActive", "'true' AS isActive")
$sqlCell.Value = $sqlText

# Leave the originally active sheet selected again.
$ws1.Activate()
